# Update the "data" sheet's time_taken column (F2:F108) with refreshed
# query timestamps, and add a new "metadata" worksheet summarising the
# panel query (data_name/data_id/data_version/... columns) placed after it.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Refresh the per-row query timestamps on the "data" sheet (col F).
# ---------------------------------------------------------------------
$timeMap = @{
    2 = "2021-10-05 14:20:51.031333"
    3 = "2021-10-05 14:20:51.031340"
    4 = "2021-10-05 14:20:51.031343"
    5 = "2021-10-05 14:20:51.031345"
    6 = "2021-10-05 14:20:51.031348"
    7 = "2021-10-05 14:20:51.031350"
    8 = "2021-10-05 14:20:51.031353"
    9 = "2021-10-05 14:20:51.031355"
    10 = "2021-10-05 14:20:51.031358"
    11 = "2021-10-05 14:20:51.031360"
    12 = "2021-10-05 14:20:51.031362"
    13 = "2021-10-05 14:20:51.031365"
    14 = "2021-10-05 14:20:51.031367"
    15 = "2021-10-05 14:20:51.031369"
    16 = "2021-10-05 14:20:51.031372"
    17 = "2021-10-05 14:20:51.031374"
    18 = "2021-10-05 14:20:51.031376"
    19 = "2021-10-05 14:20:51.031379"
    20 = "2021-10-05 14:20:51.031381"
    21 = "2021-10-05 14:20:51.031383"
    22 = "2021-10-05 14:20:51.031386"
    23 = "2021-10-05 14:20:51.031388"
    24 = "2021-10-05 14:20:51.031391"
    25 = "2021-10-05 14:20:51.031393"
    26 = "2021-10-05 14:20:51.031396"
    27 = "2021-10-05 14:20:51.031398"
    28 = "2021-10-05 14:20:51.031400"
    29 = "2021-10-05 14:20:51.031403"
    30 = "2021-10-05 14:20:51.031405"
    31 = "2021-10-05 14:20:51.031407"
    32 = "2021-10-05 14:20:51.031410"
    33 = "2021-10-05 14:20:51.031412"
    34 = "2021-10-05 14:20:51.031415"
    35 = "2021-10-05 14:20:51.031418"
    36 = "2021-10-05 14:20:51.031420"
    37 = "2021-10-05 14:20:51.031422"
    38 = "2021-10-05 14:20:51.031425"
    39 = "2021-10-05 14:20:51.031427"
    40 = "2021-10-05 14:20:51.031430"
    41 = "2021-10-05 14:20:51.031432"
    42 = "2021-10-05 14:20:51.031435"
    43 = "2021-10-05 14:20:51.031437"
    44 = "2021-10-05 14:20:51.031439"
    45 = "2021-10-05 14:20:51.031442"
    46 = "2021-10-05 14:20:51.031444"
    47 = "2021-10-05 14:20:51.031446"
    48 = "2021-10-05 14:20:51.031449"
    49 = "2021-10-05 14:20:51.031451"
    50 = "2021-10-05 14:20:51.031454"
    51 = "2021-10-05 14:20:51.031456"
    52 = "2021-10-05 14:20:51.031459"
    53 = "2021-10-05 14:20:51.031461"
    54 = "2021-10-05 14:20:51.031464"
    55 = "2021-10-05 14:20:51.031467"
    56 = "2021-10-05 14:20:51.031469"
    57 = "2021-10-05 14:20:51.031471"
    58 = "2021-10-05 14:20:51.031474"
    59 = "2021-10-05 14:20:51.031476"
    60 = "2021-10-05 14:20:51.031478"
    61 = "2021-10-05 14:20:51.031480"
    62 = "2021-10-05 14:20:51.031483"
    63 = "2021-10-05 14:20:51.031485"
    64 = "2021-10-05 14:20:51.031487"
    65 = "2021-10-05 14:20:51.031489"
    66 = "2021-10-05 14:20:51.031492"
    67 = "2021-10-05 14:20:51.031495"
    68 = "2021-10-05 14:20:51.031497"
    69 = "2021-10-05 14:20:51.031499"
    70 = "2021-10-05 14:20:51.031502"
    71 = "2021-10-05 14:20:51.031504"
    72 = "2021-10-05 14:20:51.031506"
    73 = "2021-10-05 14:20:51.031508"
    74 = "2021-10-05 14:20:51.031511"
    75 = "2021-10-05 14:20:51.031513"
    76 = "2021-10-05 14:20:51.031515"
    77 = "2021-10-05 14:20:51.031517"
    78 = "2021-10-05 14:20:51.031521"
    79 = "2021-10-05 14:20:51.031524"
    80 = "2021-10-05 14:20:51.031526"
    81 = "2021-10-05 14:20:51.031529"
    82 = "2021-10-05 14:20:51.031531"
    83 = "2021-10-05 14:20:51.031533"
    84 = "2021-10-05 14:20:51.031535"
    85 = "2021-10-05 14:20:51.031538"
    86 = "2021-10-05 14:20:51.031540"
    87 = "2021-10-05 14:20:51.031542"
    88 = "2021-10-05 14:20:51.031544"
    89 = "2021-10-05 14:20:51.031547"
    90 = "2021-10-05 14:20:51.031549"
    91 = "2021-10-05 14:20:51.031551"
    92 = "2021-10-05 14:20:51.031553"
    93 = "2021-10-05 14:20:51.031556"
    94 = "2021-10-05 14:20:51.031559"
    95 = "2021-10-05 14:20:51.031562"
    96 = "2021-10-05 14:20:51.031564"
    97 = "2021-10-05 14:20:51.031566"
    98 = "2021-10-05 14:20:51.031568"
    99 = "2021-10-05 14:20:51.031571"
    100 = "2021-10-05 14:20:51.031573"
    101 = "2021-10-05 14:20:51.031575"
    102 = "2021-10-05 14:20:51.031577"
    103 = "2021-10-05 14:20:51.031580"
    104 = "2021-10-05 14:20:51.031582"
    105 = "2021-10-05 14:20:51.031584"
    106 = "2021-10-05 14:20:51.031586"
    107 = "2021-10-05 14:20:51.031589"
    108 = "2021-10-05 14:20:51.031591"
}

foreach ($row in $timeMap.Keys) {
    $dataSheet.Cells.Item($row, 6).Value = $timeMap[$row]
}

# ---------------------------------------------------------------------
# 2. Add the new "metadata" sheet right after "data".
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$metaSheet = $wb.Worksheets.Add($null, $lastSheet)
$metaSheet.Name = "metadata"

# Reuse the bold/bordered header style already used on the "data" sheet
# (row 1 header cells + the column-A index cells) instead of defining a
# brand-new style in styles.xml.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Hydrocephalus"
$metaSheet.Range("C2").Value = 179

# "2.116" must stay a text value (matches the source export), not be
# coerced to the number 2.116 - force text via NumberFormat, assign, then
# drop back to the default (unstyled) format so no stray style is left
# behind on the cell.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "2.116"
$metaSheet.Range("D2").ClearFormats()

$metaSheet.Range("E2").Value = "2021-08-17T14:06:08.849769Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:20:51.028285"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/179/?format=json"

# Keep "data" as the active/selected sheet (matches the original file).
$dataSheet.Activate()
